$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03887266666666667
$ws.Range("H2").Value = 0.116618
$ws.Range("M2").Value = 5.273684
$ws.Range("N2").Value = 15.821052
$ws.Range("O2").Value = 0.0510821201937383
$ws.Range("P2").Value = 0.0510821201937383
$ws.Range("Q2").Value = 0.2050021602373333
$ws.Range("R2").Value = 1.845019442136
$ws.Range("S2").Value = 0.0510821201937383
$ws.Range("T2").Value = 0.0510821201937383

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03887266666666667
$ws.Range("H3").Value = 0.116618
$ws.Range("O3").Value = 0.5598845502029881
$ws.Range("P3").Value = 0.5598845502029881
$ws.Range("Q3").Value = 2.246922050999555
$ws.Range("R3").Value = 20.222298458996
$ws.Range("S3").Value = 0.5598845502029881
$ws.Range("T3").Value = 0.5598845502029881

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03887266666666667
$ws.Range("H4").Value = 0.116618
$ws.Range("M4").Value = 32.95839133333334
$ws.Range("N4").Value = 98.87517400000002
$ws.Range("O4").Value = 0.3192425840231603
$ws.Range("P4").Value = 0.3192425840231604
$ws.Range("Q4").Value = 1.281180560170222
$ws.Range("R4").Value = 11.530625041532
$ws.Range("S4").Value = 0.3192425840231603
$ws.Range("T4").Value = 0.3192425840231604

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.03887266666666667
$ws.Range("H5").Value = 0.116618
$ws.Range("M5").Value = 7.205150000000001
$ws.Range("N5").Value = 21.61545
$ws.Range("O5").Value = 0.06979074558011317
$ws.Range("P5").Value = 0.06979074558011318
$ws.Range("Q5").Value = 0.2800833942333333
$ws.Range("R5").Value = 2.5207505481
$ws.Range("S5").Value = 0.06979074558011317
$ws.Range("T5").Value = 0.06979074558011318
